# Homework 3 changes:
#  1. Footer/header "Updated automatically" date field on the slide master
#     and every slide layout moves from 9/25/2018 -> 9/30/2018.
#  2. The "Arc 45" autoshape on slide 1 is rotated / resized / reshaped.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached text of the "datetimeFigureOut" date placeholder
#    field everywhere it appears: the slide master and all custom
#    (slide) layouts.
# ---------------------------------------------------------------------
$newDate = "9/30/2018"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Re-transform the "Arc 45" shape on slide 1.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$arc = $slide.Shapes.Item("Arc 45")

$arc.Rotation = 81.13221666666666
$arc.Left = 439.2676377952756
$arc.Top = 293.31173228346455
$arc.Width = 139.2010300
$arc.Height = 162.7505600

# Adjustments must be written highest-index-first: writing a lower index
# after a higher one keeps both values, the reverse order clobbers adj1.
$arc.Adjustments.Item(2) = 1.69469
$arc.Adjustments.Item(1) = 179.07076
